$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (interested count) in column F
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 55
$wsExpo.Range("F3").Value = 96
$wsExpo.Range("F4").Value = 7282
$wsExpo.Range("F5").Value = 268
$wsExpo.Range("F6").Value = 425
$wsExpo.Range("F7").Value = 3775
$wsExpo.Range("F10").Value = 269
$wsExpo.Range("F11").Value = 609
$wsExpo.Range("F12").Value = 102

# Sheet "全部类型" (All types) - same underlying events, update the matching rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 55
$wsAll.Range("F3").Value = 96
$wsAll.Range("F5").Value = 7282
$wsAll.Range("F7").Value = 268
$wsAll.Range("F8").Value = 425
$wsAll.Range("F9").Value = 3775
$wsAll.Range("F12").Value = 269
$wsAll.Range("F13").Value = 609
$wsAll.Range("F14").Value = 102
